$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the V1/V2/V3 Time headers to the new GPU-solver labels ---
$ws.Range("H2").Value = "CPU"
$ws.Range("I2").Value = "GPU w/CPU reduc"
$ws.Range("J2").Value = "GPU w/GPU reduc"

# --- Fill in the new benchmark row (row 3) with the GPU timings ---
# Pull formatting (borders, no fill) from cells that already carry the
# "clean" (unshaded) look we want for H3:J3, then write the values.
$ws.Range("A3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("J3").PasteSpecial(-4122)

$ws.Range("H3").Value = 2.957
$ws.Range("I3").Value = 2.261
$ws.Range("J3").Value = 1.411

# --- Give the bottom row (row 5, H:J) the shaded look used elsewhere ---
$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("I4").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("J4").Copy()
$ws.Range("J5").PasteSpecial(-4122)

# --- Column widths for the new I/J columns ---
$ws.Columns.Item(9).ColumnWidth = 14.17
$ws.Columns.Item(10).ColumnWidth = 14.33

Write-Output "done"
